# ---------------------------------------------------------------------------
# update ndn consumer producer api analysis
#
# 1. Insert six new paragraphs at the very top of the document: a title
#    line, two note lines, an error line that carries the "_GoBack"
#    bookmark, and two blank paragraphs.
# 2. Merge the two runs of the "AddToWallet ..." paragraph into a single
#    run (the bookmark that used to sit at the end of that paragraph moved
#    up to the newly inserted "403forbidden" paragraph instead).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. insert the new leading paragraphs -----------------------------------

$insertionPoint = $d.Paragraphs.First.Range.Duplicate
$insertionPoint.Collapse(1)

# "Z" is a throw-away marker character. It gives us a non-paragraph-mark
# position to anchor the relocated bookmark on (and to later trim down to an
# empty paragraph) without ever handing a zero-length Range that sits exactly
# on a paragraph mark to InsertBefore/Bookmarks.Add. The marker is deleted
# again immediately after it has served its purpose.
$newText = "搭建runConsensus 遇到的错误`r" + `
           "1 sudo vi /etc/resolv.comf  要把dns服务器提到第一位`r" + `
           "2 出现各种连接错误的情况：`r" + `
           "403forbidden    bitcoind没有开启网段接收Z`r" + `
           "Z`r" + `
           "Z`r"

$insertionPoint.InsertBefore($newText)

# Drop the "_GoBack" bookmark right after the "403forbidden..." text.
$errParagraph = $d.Paragraphs(4).Range
$markerPos = $errParagraph.End - 2
$bookmarkRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
$d.Range($markerPos, $markerPos + 1).Delete()

# Strip the marker characters out of the two trailing blank paragraphs.
$blank1 = $d.Paragraphs(5).Range
$d.Range($blank1.Start, $blank1.Start + 1).Delete()

$blank2 = $d.Paragraphs(6).Range
$d.Range($blank2.Start, $blank2.Start + 1).Delete()

# --- 2. merge the AddToWallet runs -------------------------------------------

# The paragraph currently holds two adjacent runs with identical formatting
# followed by the old "_GoBack" bookmark:
#   "...wallet.dat" + "，只有发送/接收金币方调用此函数" + _GoBack
# Re-typing the whole (unchanged) sentence over itself collapses the two runs
# back into one run and consumes the now-stale bookmark in the process.
$mergedText = "AddToWallet 函数用来检测是否插入交易，fInsertedNew 判断是否插入到mapWallet缓存  fUpdated判断是否更新钱包到 wallet.dat，只有发送/接收金币方调用此函数"
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)

# Safety net: if the old bookmark somehow survived the re-typing above, make
# sure only it (and not the one we just planted) gets removed.
if ($d.Bookmarks.Exists("_GoBack")) {
    $current = $d.Bookmarks("_GoBack")
    if ($current.Range.Start -ne $bookmarkRange.Start) {
        $current.Delete()
    }
}
